{"js": "// Update the date heading paragraph (first paragraph of the body).\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst heading = paras.items[0];\nheading.load(\"text\");\nawait context.sync();\n\nif (heading.text.trim() === \"2024-09-21 Saturday\") {\n  heading.insertText(\"2024-09-22 Sunday\", Word.InsertLocation.replace);\n}\n\n// Update the worksheet's answer cells. Each entry is\n// [rowIndex, columnIndex, expectedOldValue, newValue] using the single\n// table's 0-based row/column addressing (matching Table.getCell).\nconst cellUpdates = [\n  [0, 0, \"19\u00f74=4, 3\", \"43\u00f72=21, 1\"],\n  [0, 1, \"12\u00f73=4, 0\", \"47\u00f78=5, 7\"],\n  [0, 2, \"36\u00f73=12, 0\", \"45\u00f72=22, 1\"],\n  [0, 3, \"84\u00f78=10, 4\", \"29\u00f75=5, 4\"],\n  [0, 4, \"94\u00f72=47, 0\", \"46\u00f72=23, 0\"],\n  [4, 0, \"44\u00f79=4, 8\", \"32\u00f72=16, 0\"],\n  [4, 1, \"86\u00f74=21, 2\", \"70\u00f72=35, 0\"],\n  [4, 2, \"38\u00f77=5, 3\", \"57\u00f79=6, 3\"],\n  [4, 3, \"92\u00f72=46, 0\", \"19\u00f75=3, 4\"],\n  [4, 4, \"14\u00f74=3, 2\", \"97\u00f79=10, 7\"],\n  [8, 0, \"21\u00f76=3, 3\", \"72\u00f72=36, 0\"],\n  [8, 1, \"78\u00f79=8, 6\", \"64\u00f78=8, 0\"],\n  [8, 2, \"20\u00f75=4, 0\", \"53\u00f75=10, 3\"],\n  [8, 3, \"86\u00f73=28, 2\", \"69\u00f73=23, 0\"],\n  [8, 4, \"91\u00f79=10, 1\", \"23\u00f77=3, 2\"],\n  [12, 0, \"83\u00f78=10, 3\", \"32\u00f74=8, 0\"],\n  [12, 1, \"20\u00f73=6, 2\", \"11\u00f77=1, 4\"],\n  [12, 2, \"17\u00f75=3, 2\", \"97\u00f74=24, 1\"],\n  [12, 3, \"55\u00f74=13, 3\", \"94\u00f75=18, 4\"],\n  [12, 4, \"79\u00f73=26, 1\", \"39\u00f79=4, 3\"],\n  [16, 0, \"85\u00f72=42, 1\", \"92\u00f75=18, 2\"],\n  [16, 1, \"97\u00f72=48, 1\", \"52\u00f73=17, 1\"],\n  [16, 2, \"32\u00f75=6, 2\", \"16\u00f75=3, 1\"],\n  [16, 3, \"18\u00f79=2, 0\", \"55\u00f73=18, 1\"],\n  [16, 4, \"19\u00f75=3, 4\", \"33\u00f78=4, 1\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst cells = cellUpdates.map(([row, col]) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\ncells.forEach((cell, i) => {\n  const [, , oldValue, newValue] = cellUpdates[i];\n  if (cell.value === oldValue) {\n    cell.value = newValue;\n  }\n});\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph of the document body).\n$heading = $d.Paragraphs.Item(1).Range\nif ($heading.Text.TrimEnd(\"`r\") -eq \"2024-09-21 Saturday\") {\n    $heading.Text = \"2024-09-22 Sunday\"\n}\n\n# Update the worksheet's answer cells. Each row is\n# (table row, table column, expected old value, new value) using Word's\n# 1-based Table.Cell(row, column) addressing.\n$updates = @(\n    @{Row = 1;  Col = 1; Old = \"19\u00f74=4, 3\";   New = \"43\u00f72=21, 1\"},\n    @{Row = 1;  Col = 2; Old = \"12\u00f73=4, 0\";   New = \"47\u00f78=5, 7\"},\n    @{Row = 1;  Col = 3; Old = \"36\u00f73=12, 0\";  New = \"45\u00f72=22, 1\"},\n    @{Row = 1;  Col = 4; Old = \"84\u00f78=10, 4\";  New = \"29\u00f75=5, 4\"},\n    @{Row = 1;  Col = 5; Old = \"94\u00f72=47, 0\";  New = \"46\u00f72=23, 0\"},\n    @{Row = 5;  Col = 1; Old = \"44\u00f79=4, 8\";   New = \"32\u00f72=16, 0\"},\n    @{Row = 5;  Col = 2; Old = \"86\u00f74=21, 2\";  New = \"70\u00f72=35, 0\"},\n    @{Row = 5;  Col = 3; Old = \"38\u00f77=5, 3\";   New = \"57\u00f79=6, 3\"},\n    @{Row = 5;  Col = 4; Old = \"92\u00f72=46, 0\";  New = \"19\u00f75=3, 4\"},\n    @{Row = 5;  Col = 5; Old = \"14\u00f74=3, 2\";   New = \"97\u00f79=10, 7\"},\n    @{Row = 9;  Col = 1; Old = \"21\u00f76=3, 3\";   New = \"72\u00f72=36, 0\"},\n    @{Row = 9;  Col = 2; Old = \"78\u00f79=8, 6\";   New = \"64\u00f78=8, 0\"},\n    @{Row = 9;  Col = 3; Old = \"20\u00f75=4, 0\";   New = \"53\u00f75=10, 3\"},\n    @{Row = 9;  Col = 4; Old = \"86\u00f73=28, 2\";  New = \"69\u00f73=23, 0\"},\n    @{Row = 9;  Col = 5; Old = \"91\u00f79=10, 1\";  New = \"23\u00f77=3, 2\"},\n    @{Row = 13; Col = 1; Old = \"83\u00f78=10, 3\";  New = \"32\u00f74=8, 0\"},\n    @{Row = 13; Col = 2; Old = \"20\u00f73=6, 2\";   New = \"11\u00f77=1, 4\"},\n    @{Row = 13; Col = 3; Old = \"17\u00f75=3, 2\";   New = \"97\u00f74=24, 1\"},\n    @{Row = 13; Col = 4; Old = \"55\u00f74=13, 3\";  New = \"94\u00f75=18, 4\"},\n    @{Row = 13; Col = 5; Old = \"79\u00f73=26, 1\";  New = \"39\u00f79=4, 3\"},\n    @{Row = 17; Col = 1; Old = \"85\u00f72=42, 1\";  New = \"92\u00f75=18, 2\"},\n    @{Row = 17; Col = 2; Old = \"97\u00f72=48, 1\";  New = \"52\u00f73=17, 1\"},\n    @{Row = 17; Col = 3; Old = \"32\u00f75=6, 2\";   New = \"16\u00f75=3, 1\"},\n    @{Row = 17; Col = 4; Old = \"18\u00f79=2, 0\";   New = \"55\u00f73=18, 1\"},\n    @{Row = 17; Col = 5; Old = \"19\u00f75=3, 4\";   New = \"33\u00f78=4, 1\"}\n)\n\n$table = $d.Tables.Item(1)\nforeach ($u in $updates) {\n    $cell = $table.Cell($u.Row, $u.Col)\n    $cellRange = $cell.Range\n    $cellText = $cellRange.Text.TrimEnd(\"`r\", \"`a\")\n    if ($cellText -eq $u.Old) {\n        $cellRange.Text = $u.New\n    }\n}\n"}
